# Apply updated odds values to Sheet1 as described by the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 3
$ws.Range("M3").Value = 1.1
$ws.Range("N3").Value = 7
$ws.Range("Q3").Value = 2.5
$ws.Range("R3").Value = 1.5

# Row 6
$ws.Range("I6").Value = 6.25
$ws.Range("AC6").Value = 9
$ws.Range("AD6").Value = 7.5
$ws.Range("AH6").Value = 13
$ws.Range("AI6").Value = 29

# Row 8
$ws.Range("M8").Value = 1.04
$ws.Range("N8").Value = 13
$ws.Range("Q8").Value = 1.83
$ws.Range("R8").Value = 2.03

# Row 9
$ws.Range("G9").Value = 1.44
$ws.Range("H9").Value = 4
$ws.Range("I9").Value = 8
$ws.Range("L9").Value = 7.5
$ws.Range("W9").Value = 5.5
$ws.Range("Y9").Value = 9
$ws.Range("AD9").Value = 8
$ws.Range("AV9").Value = 81
$ws.Range("AW9").Value = 8.5

# Row 10
$ws.Range("G10").Value = 1.73
$ws.Range("H10").Value = 3.3
$ws.Range("I10").Value = 5.25
$ws.Range("U10").Value = 2
$ws.Range("V10").Value = 1.73
$ws.Range("X10").Value = 7.5
$ws.Range("AC10").Value = 8
$ws.Range("AJ10").Value = 17
$ws.Range("AN10").Value = 3.6
$ws.Range("AU10").Value = 9

# Row 12
$ws.Range("N12").Value = 5
